# Generate Report for Handoff
# Updates localization status report: status transitions, handoff timestamps,
# and error details for out-of-date handback files.

$wb = $excel.ActiveWorkbook

$errorDetail1 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc6022a08dcdde835deef71a5010e0de3d6cb3d8/e2e/1c7e79d1-ef9b-47d4-ad9e-73355539705a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94275538884a7f2c5b814f0bb4e9c04ab13fc8b6/e2e/1c7e79d1-ef9b-47d4-ad9e-73355539705a.md."
$errorDetail2 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bc6022a08dcdde835deef71a5010e0de3d6cb3d8/e2e/9c7ecf02-c54f-45b0-bd12-16d5d6aa1bf8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/94275538884a7f2c5b814f0bb4e9c04ab13fc8b6/e2e/9c7ecf02-c54f-45b0-bd12-16d5d6aa1bf8.md."

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-10-26 08:16:25"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-26 08:16:25"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(16).ColumnWidth = 40

$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2016-10-26 08:16:12"
$wsZhCn.Range("P2").Value = $errorDetail1

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-10-26 08:16:12"
$wsZhCn.Range("P3").Value = $errorDetail2

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(16).ColumnWidth = 40

$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("H2").Value = "2016-10-26 08:16:25"
$wsDeDe.Range("P2").Value = $errorDetail1

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-10-26 08:16:25"
$wsDeDe.Range("P3").Value = $errorDetail2
